$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "51.629.07"
Set-TextValue "E2" "  +1.12%  "
Set-TextValue "D3" "3.021.00"
Set-TextValue "E3" "  +2.12%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "379.19"
Set-TextValue "E5" "  +0.04%  "
Set-TextValue "D6" "102.59"
Set-TextValue "E6" "  +0.00%  "
Set-TextValue "D7" "0.547"
Set-TextValue "E7" "  +0.45%  "
Set-TextValue "E8" "  +0.01%  "
Set-TextValue "D9" "0.590"
Set-TextValue "E9" "  +0.71%  "
Set-TextValue "D10" "36.73"
Set-TextValue "E10" "  +0.74%  "
Set-TextValue "E11" "  -0.16%  "
Set-TextValue "D12" "0.0863"
Set-TextValue "D13" "3.502.41"
Set-TextValue "E13" "  +1.90%  "
Set-TextValue "D14" "18.43"
Set-TextValue "E14" "  +0.04%  "
Set-TextValue "D15" "7.72"
Set-TextValue "D16" "3.023.14"
Set-TextValue "E16" "  +2.23%  "
Set-TextValue "D17" "0.974"
Set-TextValue "E17" "  -3.92%  "
Set-TextValue "D18" "10.62"
Set-TextValue "E18" "  -14.69%  "
Set-TextValue "D19" "51.627.88"
Set-TextValue "E19" "  +1.03%  "
Set-TextValue "E20" "  +0.61%  "
Set-TextValue "E21" "  +0.11%  "
Set-TextValue "E22" "  +0.85%  "
Set-TextValue "D23" "70.00"
Set-TextValue "E23" "  +0.49%  "
Set-TextValue "D24" "267.48"
Set-TextValue "E24" "  -0.03%  "
Set-TextValue "E25" "  -6.15%  "
Set-TextValue "D26" "8.32"
Set-TextValue "E26" "  +4.11%  "
Set-TextValue "E27" "  +7.65%  "
Set-TextValue "D28" "0.173"
Set-TextValue "E28" "  +4.45%  "
Set-TextValue "E29" "  +0.04%  "
Set-TextValue "D30" "26.19"
Set-TextValue "E30" "  +1.39%  "
Set-TextValue "E31" "  +0.31%  "
Set-TextValue "D32" "10.26"
Set-TextValue "E32" "  -2.38%  "
Set-TextValue "D33" "2.12"
Set-TextValue "E33" "  +2.82%  "
Set-TextValue "D34" "50.56"
Set-TextValue "E34" "  -0.43%  "
Set-TextValue "D35" "33.84"
Set-TextValue "E35" "  -0.97%  "
Set-TextValue "E36" "  +3.30%  "
Set-TextValue "E38" "  +2.21%  "
Set-TextValue "D39" "0.292"
Set-TextValue "E39" "  +13.47%  "
Set-TextValue "D40" "16.90"
Set-TextValue "E40" "  +0.69%  "
Set-TextValue "E41" "  +1.45%  "
Set-TextValue "B42" "Monero"
Set-TextValue "C42" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D42" "127.67"
Set-TextValue "E42" "  +7.49%  "
Set-TextValue "B43" "Stellar"
Set-TextValue "C43" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D43" "0.116"
Set-TextValue "E43" "  -0.66%  "
Set-TextValue "E44" "  +1.61%  "
Set-TextValue "D45" "3.79"
Set-TextValue "E45" "  +5.61%  "
Set-TextValue "E46" "  -0.87%  "
Set-TextValue "E47" "  +2.71%  "
Set-TextValue "E48" "  +2.37%  "
Set-TextValue "D49" "2.025.81"
Set-TextValue "E49" "  -0.86%  "
Set-TextValue "D50" "3.319.73"
Set-TextValue "E50" "  +2.09%  "
Set-TextValue "E51" "  -1.58%  "
